$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.730.53'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '1.594.74'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.502'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.50%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.62%  '
$ws.Range("E9").Value = '  -2.06%  '
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").Value = '1.819.79'
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").Value = '1.591.57'
$ws.Range("E13").Value = '  -2.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.533'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.78%  '
$ws.Range("D16").Value = '27.692.88'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.59%  '
$ws.Range("D19").Value = '0.0₃0697'
$ws.Range("E19").Value = '  -2.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.33%  '
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("E24").Value = '  -3.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.63%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.106'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.69%  '
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0472'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.97%  '
$ws.Range("D33").Value = '1.376.99'
$ws.Range("E33").Value = '  -2.93%  '
$ws.Range("E34").Value = '  -4.43%  '
$ws.Range("E35").Value = '  -4.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.970'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.16%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  -1.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.538'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.831'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.59%  '
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.973'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("E44").Value = '  +2.48%  '
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("E46").Value = '  -5.11%  '
$ws.Range("D47").Value = '1.730.59'
$ws.Range("E47").Value = '  -2.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.98'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.03%  '
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0967'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0496'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.28%  '
